$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.347.05'
$ws.Range("E2").Value = '  -4.68%  '
$ws.Range("D3").Value = '3.255.73'
$ws.Range("E3").Value = '  -7.58%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '590.40'
$ws.Range("E5").Value = '  -5.24%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.67'
$ws.Range("E6").Value = '  -12.52%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '3.245.58'
$ws.Range("E8").Value = '  -7.77%  '
$ws.Range("E9").Value = '  -11.02%  '
$ws.Range("E10").Value = '  -13.33%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.82'
$ws.Range("E11").Value = '  -3.01%  '
$ws.Range("E12").Value = '  -12.82%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.38'
$ws.Range("E13").Value = '  -17.17%  '
$ws.Range("E14").Value = '  -12.11%  '
$ws.Range("D15").Value = '3.771.31'
$ws.Range("E15").Value = '  -7.74%  '
$ws.Range("D16").Value = '67.292.84'
$ws.Range("E16").Value = '  -4.92%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '549.13'
$ws.Range("E17").Value = '  -9.96%  '
$ws.Range("D18").Value = '3.252.45'
$ws.Range("E18").Value = '  -7.53%  '
$ws.Range("E19").Value = '  -13.44%  '
$ws.Range("E20").Value = '  -5.92%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.20'
$ws.Range("E21").Value = '  -14.14%  '
$ws.Range("E22").Value = '  -13.10%  '
$ws.Range("E23").Value = '  -14.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.67'
$ws.Range("E24").Value = '  -12.90%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.57'
$ws.Range("E25").Value = '  -12.97%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("E27").Value = '  -14.27%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '29.46'
$ws.Range("E28").Value = '  -12.51%  '
$ws.Range("E29").Value = '  -10.93%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.13'
$ws.Range("E30").Value = '  -16.82%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.67'
$ws.Range("E31").Value = '  -11.14%  '
$ws.Range("E32").Value = '  -12.59%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '546.72'
$ws.Range("E33").Value = '  -14.95%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.67'
$ws.Range("E34").Value = '  -17.61%  '
$ws.Range("E35").Value = '  -15.37%  '
$ws.Range("E36").Value = '  -0.21%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0450'
$ws.Range("E37").Value = '  -4.57%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '53.73'
$ws.Range("E38").Value = '  -5.58%  '
$ws.Range("E39").Value = '  -13.99%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '9.21'
$ws.Range("E40").Value = '  -14.54%  '
$ws.Range("E41").Value = '  -11.62%  '
$ws.Range("D42").Value = '2.932.20'
$ws.Range("E42").Value = '  -12.34%  '
$ws.Range("E43").Value = '  -23.15%  '
$ws.Range("E44").Value = '  -15.44%  '
$ws.Range("E45").Value = '  -19.14%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '26.58'
$ws.Range("E46").Value = '  -16.55%  '
$ws.Range("E47").Value = '  -15.16%  '
$ws.Range("E48").Value = '  -0.12%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '127.84'
$ws.Range("E49").Value = '  -4.41%  '
$ws.Range("E50").Value = '  -20.18%  '
$ws.Range("E51").Value = '  -12.19%  '
